$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F4").Value = "UserAuthentication newUser(String nickname)"

$ws.Range("F17").Value = "WordBank selectWords()"
$ws.Range("F18").Value = "WordBank selectWords()"
$ws.Range("F19").Value = "WordBank selectWords()"
$ws.Range("F20").Value = "WordBank selectWords()"

$ws.Range("D21").Value = "The words to find will be horizontal, vertical down, vertical up, diagonal down, or diagonal up"

$ws.Range("F26").Value = "GameSession getWordPositions()"
$ws.Range("F27").Value = "Index.html highlightButton(row, column, idx)"
$ws.Range("F31").Value = "Index.html highlightButton(row, column, idx)"
$ws.Range("F32").Value = "Index.html highlightButton(row, column, idx)"

$ws.Range("F34").Value = "GameSession charSelected(int, int)"
$ws.Range("F35").Value = "GameSession charSelected(int, int)"
$ws.Range("F36").Value = "GameSession charSelected(int, int)"

$ws.Range("D40").Value = "Each possible orientation of words will be at least 15% of the total word count"

$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

$ws.Range("D41").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
